$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: A3 now carries the "AD.SEC.002.FON.01" context value (previously RO.ACT.001)
$ws.Range("A3").Value = "AD.SEC.002.FON.01"

# Rows 4 & 5: drop the now-redundant context values in column A
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()

# Row 6: drop its lone context value entirely (row collapses once empty)
$ws.Range("A6").ClearContents()

# Re-add the context info on a new row (7), this time in column D,
# matching the "header" text style used elsewhere (e.g. D2/A2)
$ws.Range("D7").Value = "AD.SEC.001.FON.01"
$ws.Range("D7").NumberFormat = "@"

# Update the remembered selection to match the new point of interest
$ws.Range("A5").Select()
